$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K10 should become a real number (was inline string "360371")
$ws.Range("K10").Value = 360371

# New row 11
$ws.Range("A11").Value = "kotak"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "+919409727736"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").Value = "L.ranpura@kotak.com"
$ws.Range("D11").Value = "L.ranpura, kotak.com, www.kotak.com"
$ws.Range("E11").Value = "Darshan Ranpura Service"
$ws.Range("F11").Value = "360 002, India Main Road"
$ws.Range("G11").Value = "kotak Kotak Mahindra Bank Darshan Ranpura Service Officer Assistant Manager Kotak Mahindra Bank Ltd_ D +91281 2812581401 Ground Floor; Chandra Cottage M+91 9409727736 Opp RMC Swimming Pool L.ranpura@kotak.com darshan. Near Kalola Children Hospital, Kothariya www.kotak.com 360 002, India Main Road, Rajkot"
$ws.Range("H11").Value = "2025-05-15 13:01:31"
$ws.Range("I11").Value = "kotak (size: 10.00), Kotak Mahindra Bank (size: 7.85), Darshan Ranpura (size: 4.52), Service Officer (size: 1.59), Assistant Manager (size: 2.96)"
$ws.Range("J11").Value = "Not Found"
$ws.Range("K11").Value = "Not Found"

# New row 12
$ws.Range("A12").Value = "kotak"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "+919409727736"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = "L.ranpura@kotak.com"
$ws.Range("D12").Value = "L.ranpura, kotak.com, www.kotak.com"
$ws.Range("E12").Value = "Darshan Ranpura Service"
$ws.Range("F12").Value = "360 002, India Main Road"
$ws.Range("G12").Value = "kotak Kotak Mahindra Bank Darshan Ranpura Service Officer Assistant Manager Kotak Mahindra Bank Ltd_ D +91281 2812581401 Ground Floor; Chandra Cottage M+91 9409727736 Opp RMC Swimming Pool L.ranpura@kotak.com darshan. Near Kalola Children Hospital, Kothariya www.kotak.com 360 002, India Main Road, Rajkot"
$ws.Range("H12").Value = "2025-05-15 13:01:39"
$ws.Range("I12").Value = "kotak (size: 10.00), Kotak Mahindra Bank (size: 7.85), Darshan Ranpura (size: 4.52), Service Officer (size: 1.59), Assistant Manager (size: 2.96)"
$ws.Range("J12").Value = "Not Found"
$ws.Range("K12").Value = "Not Found"
